$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"

# B3 holds a long digit string that must stay text (not be coerced to a
# number). Force a Text number format, assign it, then restore the
# original "General" formatting (copied from a sibling cell that already
# uses style index 8) so only the value - not the style - changes.
$cardCell = $ws.Range("B3")
$cardCell.NumberFormat = "@"
$cardCell.Value2 = "2570314725427075"
$ws.Range("B2").Copy()
$cardCell.PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 12.10.2023"

# --- Row 6 ---
$ws.Range("B6").Value = "13.10."
$ws.Range("C6").Value = "14.10."
$ws.Range("D6").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E6").Value = "24,69-"

# --- Row 7 ---
$ws.Range("B7").Value = "15.10."
$ws.Range("C7").Value = "16.10."
$ws.Range("D7").Value = "RECHNUNG VODAFONE GMBH 48255531"
$ws.Range("E7").Value = "39,00-"

# --- Row 8 ---
$ws.Range("B8").Value = "19.10."
$ws.Range("C8").Value = "20.10."
$ws.Range("D8").Value = "KARTENZ./19.10 EDEKA RO"
$ws.Range("E8").Value = "118,14-"

# --- Row 9 ---
$ws.Range("B9").Value = "21.10."
$ws.Range("C9").Value = "22.10."
$ws.Range("D9").Value = "PAYPAL PFTUSH"
$ws.Range("E9").Value = "37,22-"

# --- Row 10 ---
$ws.Range("B10").Value = "23.10."
$ws.Range("C10").Value = "24.10."
$ws.Range("D10").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 84702767"
$ws.Range("E10").Value = "85,81-"

# --- Row 11 (was blank, now populated with a new transaction) ---
$ws.Range("B11").Value = "27.10."
$ws.Range("C11").Value = "28.10."
$ws.Range("D11").Value = "BEITRAG Allianz SE K-35506483"
$ws.Range("E11").Value = "54,08-"
# The amount cell picks up the same right-aligned number style used by
# the other amount cells in column E (copied from E10, style index 17).
$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)  # xlPasteFormats

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 31.10.2023"
$ws.Range("E12").Value = "358,94-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 09.11.2023"

$excel.CutCopyMode = 0
